# manager can get presence report
# Fix a handful of incorrect arrival/month values in the presence report:
#  - row 2 (index 1, tair hadad):  arrival date-time month was Dec, should be Nov; month col -> 11
#  - row 3 (index 2, asaf rdt):    month col "12" (text) -> 11 (number)
#  - row 5 (index 4, omri sss):    arrival date-time month was Dec, should be Oct; month col -> 10
#  - row 8 (index 7, yoni machluf):arrival date-time month was Dec, should be Oct; month col -> 10

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("presence")

$ws.Range("D2").Value = "Sun, 16 Nov 2018 18:08:16"
$ws.Range("E2").Value = 11

$ws.Range("E3").Value = 11

$ws.Range("D5").Value = "Sun, 16 Oct 2018 21:08:16"
$ws.Range("E5").Value = 10

$ws.Range("D8").Value = "Wed, 19 Oct 2018 11:12:36"
$ws.Range("E8").Value = 10

$ws.Range("D15").Select() | Out-Null
